$d = $word.ActiveDocument

# remove old _GoBack bookmark (was near ypms/Demo) FIRST
$d.Bookmarks("_GoBack").Delete()

# cyxz cell: set pPr/rPr theme fonts + insert myTest11 run + bookmarks
$b = $d.Bookmarks("cyxz")
$para = $b.Range.Paragraphs(1)
$prng = $para.Range
$frag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="108"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="仿宋_GB2312" w:hAnsiTheme="minorHAnsi" w:cs="仿宋_GB2312"/></w:rPr></w:pPr><w:bookmarkStart w:id="900" w:name="cyxz"/><w:bookmarkEnd w:id="900"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="仿宋_GB2312" w:hAnsiTheme="minorHAnsi" w:cs="仿宋_GB2312"/></w:rPr><w:t>myTest11</w:t></w:r><w:bookmarkStart w:id="901" w:name="_GoBack"/><w:bookmarkEnd w:id="901"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$prng.InsertXML($frag)

Write-Host "done"
